# APTAutomationProject / APT_MCS_CustomerUser.xlsx - "All Modules- Merged Code"
#
# 1. Insert a new worksheet "unique column names" between "customerUser" and
#    "SupplyService", listing the columns that must be null/unique.
# 2. Swap the Run flag (Yes/No) + MainDomain value on two rows of customerUser.
# 3. Swap the Run flag (Yes/No) on SupplyService row 5 and move the active
#    selection there to C5.

$wb = $excel.ActiveWorkbook

$customerUser   = $wb.Worksheets.Item("customerUser")

# --- 1. New "unique column names" sheet, inserted right after customerUser ---
$uniqueCols = $wb.Worksheets.Add($null, $customerUser)
$uniqueCols.Name = "unique column names"

# Re-fetch SupplyService by name: inserting a sheet shifts worksheet
# positions, and a reference grabbed beforehand would now resolve to
# whatever sheet occupies that old position (the new sheet itself).
$supplyService  = $wb.Worksheets.Item("SupplyService")

# Give A1 the same "section header" look used elsewhere in the workbook
# (e.g. customerUser!A1) by copying its format, then overwrite the text.
$customerUser.Range("A1").Copy($uniqueCols.Range("A1"))
$uniqueCols.Range("A1").Value = "Unique Columns"

$uniqueCols.Range("A2").Value = "Name"

$uniqueCols.Range("A3").Value = "MainDomain"
$uniqueCols.Range("B3").Value = "should be null or unique"

$uniqueCols.Range("A4").Value = "OCN"

$uniqueCols.Range("A5").Value = "editCustomerName"
$uniqueCols.Range("B5").Value = "should be null or unique"

$uniqueCols.Range("A6").Value = "editMainDomain"
$uniqueCols.Range("B6").Value = "should be null or unique"

$uniqueCols.Range("A7").Value = "editOCN"
$uniqueCols.Range("B7").Value = "should be null or unique"

$uniqueCols.Columns.Item(1).ColumnWidth = 22.75
$uniqueCols.Columns.Item(2).ColumnWidth = 27.75

# --- 2. customerUser data tweaks ---
$customerUser.Range("A3").Value = "No"
$customerUser.Range("A7").Value = "Yes"
$customerUser.Range("D7").Value = "Null"

# --- 3. SupplyService data tweak + selection ---
$supplyService.Range("A5").Value = "Yes"
$supplyService.Activate()
$supplyService.Range("C5").Select()
